$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.323.02'
$ws.Range('E2').Value = '  +6.51%  '
$ws.Range('D3').Value = '2.439.68'
$ws.Range('E3').Value = '  +6.33%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '565.76'
$ws.Range('E5').Value = '  +5.12%  '
$ws.Range('D6').Value = '142.87'
$ws.Range('E6').Value = '  +11.95%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  +4.07%  '
$ws.Range('D9').Value = '2.439.33'
$ws.Range('E9').Value = '  +6.40%  '
$ws.Range('E10').Value = '  +5.17%  '
$ws.Range('D11').Value = '5.76'
$ws.Range('E11').Value = '  +5.48%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  +7.44%  '
$ws.Range('D14').Value = '26.40'
$ws.Range('E14').Value = '  +15.15%  '
$ws.Range('D15').Value = '2.870.73'
$ws.Range('E15').Value = '  +6.11%  '
$ws.Range('D16').Value = '63.133.75'
$ws.Range('E16').Value = '  +6.43%  '
$ws.Range('E17').Value = '  +9.60%  '
$ws.Range('D18').Value = '2.434.41'
$ws.Range('E18').Value = '  +5.15%  '
$ws.Range('D19').Value = '11.22'
$ws.Range('E19').Value = '  +8.79%  '
$ws.Range('D20').Value = '340.03'
$ws.Range('E20').Value = '  +10.58%  '
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +7.43%  '
$ws.Range('E22').Value = '  +4.80%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = '65.38'
$ws.Range('E24').Value = '  +4.34%  '
$ws.Range('E25').Value = '  +3.82%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  +15.09%  '
$ws.Range('D28').Value = '8.17'
$ws.Range('E28').Value = '  +6.92%  '
$ws.Range('E29').Value = '  +13.86%  '
$ws.Range('D30').Value = '6.70'
$ws.Range('E30').Value = '  +16.95%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.83'
$ws.Range('E31').Value = '  +8.06%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0790'
$ws.Range('E32').Value = '  +11.91%  '
$ws.Range('D33').Value = '174.58'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('E34').Value = '  +12.76%  '
$ws.Range('E35').Value = '  +6.75%  '
$ws.Range('D36').Value = '18.74'
$ws.Range('E36').Value = '  +6.49%  '
$ws.Range('D37').Value = '373.02'
$ws.Range('E37').Value = '  +21.27%  '
$ws.Range('D38').Value = '4.49'
$ws.Range('E38').Value = '  +13.80%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('E41').Value = '  +14.58%  '
$ws.Range('D42').Value = '40.43'
$ws.Range('E42').Value = '  +7.73%  '
$ws.Range('D43').Value = '149.58'
$ws.Range('E44').Value = '  +9.70%  '
$ws.Range('D45').Value = '20.78'
$ws.Range('E45').Value = '  +13.51%  '
$ws.Range('D46').Value = '0.596'
$ws.Range('E46').Value = '  +5.77%  '
$ws.Range('D47').Value = '0.0961'
$ws.Range('E47').Value = '  +3.23%  '
$ws.Range('D48').Value = '0.0521'
$ws.Range('E48').Value = '  +7.25%  '
$ws.Range('E49').Value = '  +7.28%  '
$ws.Range('D50').Value = '17.91'
$ws.Range('E50').Value = '  +8.67%  '
$ws.Range('D51').Value = '0.0₆0224'
$ws.Range('E51').Value = '  +3.81%  '
